$wb = $excel.ActiveWorkbook

# 1. Change the workbook's default (Normal) font from Calibri to Arial
$normalStyle = $wb.Styles.Item("Normal")
$normalStyle.Font.Name = "Arial"

# 2. Populate sheet "13" (3rd worksheet) with the new names
$ws13 = $wb.Worksheets.Item(3)
$ws13.Range("A1").Value2 = "Daniel Magnezi"
$ws13.Range("A2").Value2 = "Dniel Mani"
$ws13.Range("A3").Value2 = "Ofri Serussi"
$ws13.Range("A4").Value2 = "Nitay Man"

# 3. Make sheet "13" the active sheet/tab, with A4 selected
$ws13.Select()
$ws13.Range("A4").Select() | Out-Null
